$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Pipeline(steps=[('scaler', None), ('selector', 'passthrough'),
                ('model',
                 AdaBoostClassifier(estimator=RandomForestClassifier(class_weight='balanced',
                                                                     max_depth=1,
                                                                     min_samples_leaf=5,
                                                                     n_estimators=10,
                                                                     random_state=42),
                                    n_estimators=10, random_state=42))])"
$ws.Range("B2").Value = 0.6666666666666666
$ws.Range("C2").Value = "{'scaler': None, 'model__n_estimators': 10, 'model__estimator__n_estimators': 10, 'model__estimator__min_samples_split': 2, 'model__estimator__min_samples_leaf': 5, 'model__estimator__max_features': 'sqrt', 'model__estimator__max_depth': 1, 'model__estimator__class_weight': 'balanced'}"
$ws.Range("D2").Value = 0.4615384615384615
$ws.Range("E2").Value = "[1 0 0 1 0 0 1 1 0 1 0 0]"
$ws.Range("F2").Value = "[0 1 1 0 1 0 1 1 1 1 1 0]"
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0.9737619047619047
$ws.Range("I2").Value = 0.006782318122629813
$ws.Range("J2").Value = 0.5766666666666665
$ws.Range("K2").Value = 0.06543294592059903

# Row 3
$ws.Range("A3").Value = "Pipeline(steps=[('scaler', StandardScaler()), ('selector', 'passthrough'),
                ('model',
                 AdaBoostClassifier(estimator=RandomForestClassifier(class_weight='balanced',
                                                                     max_depth=5,
                                                                     max_features='log2',
                                                                     min_samples_split=4,
                                                                     n_estimators=50,
                                                                     random_state=42),
                                    random_state=42))])"
$ws.Range("B3").Value = 0.5809523809523809
$ws.Range("C3").Value = "{'scaler': StandardScaler(), 'model__n_estimators': 50, 'model__estimator__n_estimators': 50, 'model__estimator__min_samples_split': 4, 'model__estimator__min_samples_leaf': 1, 'model__estimator__max_features': 'log2', 'model__estimator__max_depth': 5, 'model__estimator__class_weight': 'balanced'}"
$ws.Range("D3").Value = 0.5714285714285715
$ws.Range("E3").Value = "[1 0 1 0 0 0 0 1 1 0 1 1]"
$ws.Range("F3").Value = "[1 1 1 1 1 0 0 0 0 1 1 1]"
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 0.9792380952380952
$ws.Range("I3").Value = 0.006813451456477024
$ws.Range("J3").Value = 0.5173333333333333
$ws.Range("K3").Value = 0.05400345214138303

# Row 4
$ws.Range("A4").Value = "Pipeline(steps=[('scaler', None), ('selector', 'passthrough'),
                ('model',
                 AdaBoostClassifier(estimator=RandomForestClassifier(max_depth=1,
                                                                     max_features='log2',
                                                                     min_samples_leaf=5,
                                                                     min_samples_split=4,
                                                                     n_estimators=50,
                                                                     random_state=42),
                                    n_estimators=5, random_state=42))])"
$ws.Range("B4").Value = 0.5904761904761904
$ws.Range("C4").Value = "{'scaler': None, 'model__n_estimators': 5, 'model__estimator__n_estimators': 50, 'model__estimator__min_samples_split': 4, 'model__estimator__min_samples_leaf': 5, 'model__estimator__max_features': 'log2', 'model__estimator__max_depth': 1, 'model__estimator__class_weight': None}"
$ws.Range("D4").Value = 0.823529411764706
$ws.Range("E4").Value = "[1 0 1 1 1 1 0 1 0 1 0 1]"
$ws.Range("F4").Value = "[1 0 1 1 1 1 0 0 1 1 1 1]"
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 0.983452380952381
$ws.Range("I4").Value = 0.005501670243188103
$ws.Range("J4").Value = 0.5319999999999999
$ws.Range("K4").Value = 0.07232316098127996
